$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.974.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.622.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.652.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.848.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("E15").Value = "  -4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.961.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0738"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.845"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.515"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.771"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.760.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "54.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.24%  "
